# The workbook has a PivotTable (on "Sheet2") built from the source table
# on "Sheet1". The source commit only fixes one row of the source table and
# moves the saved viewport - it does NOT refresh the PivotTable - so switch
# to manual calculation first to keep the pivot's cached output untouched
# while we edit the source data.
$excel.Calculation = -4135   # xlCalculationManual

$wb = $excel.ActiveWorkbook

# "Sheet1" holds the photo/route source table. Row 44 is the Marina Bay /
# Helix Bridge entry - fix its photo credit to the new source "MyWoWo".
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D44").Value = "MyWoWo"

# Move the saved viewport/selection down to the newly added Kent Ridge map
# rows near the bottom of the table.
$ws.Range("D47").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
